# Round of fixes to the district councils dataset (3 errors remaining).
#
# The sheet is sorted by the "weighted_total" column (N, col 14) descending.
# Three councils had an incorrect input score; fixing the input also shifts
# their weighted_total, which in turn moves them to a new position in the
# sort order. So for each fix we: locate the council's row by its
# official-name, correct the offending score, recompute weighted_total from
# the same weights Excel used for the rest of the sheet, and finally
# re-sort the whole data range so the table stays ordered by
# weighted_total descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count()

# weighted_total = 0.15*s1_gov + 0.15*s2_m&a + 0.15*s3_c&a + 0.15*s4_coms
#                + 0.10*s5_mset + 0.05*s6_cb + 0.10*s7_dsi + 0.10*s8_est + 0.05*s9_ee
# i.e. columns E..M (5..13), weighted_total lands in column N (14).
$weights = @(0.15, 0.15, 0.15, 0.15, 0.10, 0.05, 0.10, 0.10, 0.05)

function Fix-CouncilScore($councilName, $col, $newValue) {
    for ($r = 2; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value()
        if ($name -eq $councilName) {
            $ws.Cells.Item($r, $col).Value = $newValue

            $total = 0
            for ($i = 0; $i -lt 9; $i++) {
                $cellVal = $ws.Cells.Item($r, 5 + $i).Value()
                $total = $total + ($cellVal * $weights[$i])
            }
            $ws.Cells.Item($r, 14).Value = $total
            return
        }
    }
}

# s9_ee (column M / 13) was overstated for Reigate and Banstead.
Fix-CouncilScore "Reigate and Banstead Borough Council" 13 0.75

# s9_ee (column M / 13) was overstated for Basingstoke and Deane.
Fix-CouncilScore "Basingstoke and Deane Borough Council" 13 0.5

# s3_c&a (column G / 7) was overstated for Oadby and Wigston.
Fix-CouncilScore "Oadby and Wigston Borough Council" 7 0.5714285714285714

# Re-sort the whole table by weighted_total (column N) descending so the
# three corrected councils land in their new rank position.
$rng = $ws.Range("A2:S" + $lastRow)
$sortKey = $ws.Range("N2:N" + $lastRow)
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, 0, 2, 0, 0)
$ws.Sort.SetRange($rng)
$ws.Sort.Header = 0
$ws.Sort.Apply()
